# This script applies a weekly update to the "Fruta, Macroferia Regional de
# Talca - Frutilla" price table: a brand-new price record is inserted at the
# top of the existing block of records (row 681), every subsequent existing
# record shifts down by one row, and the record that used to be last (row
# 774) is re-written, unchanged, as the new last row (775).
#
# Columns D,L,M,N,O,P,Q,R,S hold the per-record "price entry" fields; columns
# A,B,C,E,F,G,H,I,J,K,T are a constant template shared by every record in the
# block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 681
$lastRow  = 774

# Column numbers for the "price entry" fields that move as a unit.
$priceCols = @(4, 12, 13, 14, 15, 16, 17, 18, 19)   # D, L, M, N, O, P, Q, R, S

# --- 1. Snapshot the current ("before") values of the price-entry columns
#        for every existing record row, keyed by row number. -----------------
$old = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $priceCols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $old[$r] = $rowVals
}

# --- 2. Shift every existing record down by one row: new row r gets the
#        values that used to live in row r-1. Work from the bottom up so we
#        never overwrite a value before it has been read (values were
#        already captured in step 1, but we still go bottom-up for safety).
for ($r = $lastRow + 1; $r -ge $firstRow + 1; $r--) {
    $src = $old[$r - 1]
    foreach ($c in $priceCols) {
        $ws.Cells.Item($r, $c).Value2 = $src[$c]
    }
}

# --- 3. New row 775 is a brand-new sheet row: copy the constant template
#        columns from the previous last row (774) so they carry over too. ---
$templateCols = @(1, 2, 3, 5, 6, 7, 8, 9, 10, 11, 20)  # A,B,C,E,F,G,H,I,J,K,T
foreach ($c in $templateCols) {
    $ws.Cells.Item($lastRow + 1, $c).Value2 = $ws.Cells.Item($lastRow, $c).Value2
}

# The "fecha" column (D) is date-formatted; the new row needs that same
# number format explicitly, since a brand-new cell has no formatting yet.
$ws.Cells.Item($lastRow + 1, 4).NumberFormat = $ws.Cells.Item($lastRow, 4).NumberFormat

# --- 4. Write the brand-new price record into the now-vacated first row
#        (681). -------------------------------------------------------------
$ws.Cells.Item($firstRow, 4).Value2  = 44984            # D - fecha
$ws.Cells.Item($firstRow, 12).Value  = "Primera"        # L - calidad
$ws.Cells.Item($firstRow, 13).Value2 = 120               # M - cantidad
$ws.Cells.Item($firstRow, 14).Value2 = 7000              # N - precio minimo
$ws.Cells.Item($firstRow, 15).Value2 = 7000              # O - precio maximo
$ws.Cells.Item($firstRow, 16).Value2 = 7000              # P - precio promedio
$ws.Cells.Item($firstRow, 17).Value  = "$/caja 7 kilos"  # Q - unidad
$ws.Cells.Item($firstRow, 18).Value  = "Región del Maule" # R - origen
$ws.Cells.Item($firstRow, 19).Value2 = 1000              # S - precio kilo
